# Update 1.3.2: Paginated loader and 2 update reference-books
#
# - Rename the single worksheet from "Sheet1" to "Материалы".
# - Move the active cell selection from D6 to E6.

$wb = $excel.ActiveWorkbook
$ws = $wb.ActiveSheet

# Rename the sheet.
$ws.Name = "Материалы"

# Update the selected/active cell on the sheet.
$ws.Range("E6").Select()
